$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added "LM cigar calculations" block below the existing summary table
# (rows 26-31): a small block re-using the CO / CO2 labels plus three
# brand-new labels (LM/LM00, ratio co2, concen co2, next).

$ws.Range("A26").Value = "LM/LM00"

$ws.Range("A27").Value = "CO"
$ws.Range("B27").Formula = "=D2/E2"

$ws.Range("A28").Value = "CO2"
$ws.Range("B28").Formula = "=D3/E3"

$ws.Range("A29").Value = "ratio co2"
$ws.Range("B29").Formula = "=100*B28"

$ws.Range("A30").Value = "concen co2"
$ws.Range("B30").Formula = "=4.1*B28"

$ws.Range("A31").Value = "next"
$ws.Range("B31").Formula = "=B27*1.2"

# Best-effort view-state nudge so the window scrolls down to the new block
# and the selection lands where the author left off (B32, just past the
# new data). Some hosts don't persist scroll position to the saved file,
# but the explicit cell selection below always does.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 4
    $win.ScrollColumn = 1
} catch {
    # Window scroll state isn't exposed in every host; ignore if unsupported.
}

$ws.Range("B32").Select()
